# Refresh the cryptos table (prices + 1h volume %) with the latest scrape,
# and fix the WrappedeETH/PEPE row ordering (rows 28-29 swapped ranks).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are stored as plain text in the source data (e.g. "68.361.24",
# "0.0000105") - mark each one as Text before writing so Excel keeps the exact
# string instead of reinterpreting it as a Double (which would drop trailing
# zeros like "5.40" -> 5.4 or switch to scientific notation).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D31", "D33", "D38", "D39", "D40", "D42", "D43", "D44", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "68.361.24"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").Value = "2.641.70"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "600.21"
$ws.Range("E5").Value = "  +1.41%  "

# Row 6
$ws.Range("D6").Value = "154.72"
$ws.Range("E6").Value = "  +2.84%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +0.40%  "

# Row 9
$ws.Range("D9").Value = "2.640.77"
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("E10").Value = "  +7.31%  "

# Row 11
$ws.Range("E11").Value = "  -0.57%  "

# Row 13
$ws.Range("E13").Value = "  +2.13%  "

# Row 14
$ws.Range("D14").Value = "28.07"
$ws.Range("E14").Value = "  +2.89%  "

# Row 15
$ws.Range("E15").Value = "  +3.66%  "

# Row 16
$ws.Range("D16").Value = "3.122.26"
$ws.Range("E16").Value = "  +1.32%  "

# Row 17
$ws.Range("D17").Value = "68.211.02"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18
$ws.Range("D18").Value = "2.631.05"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
$ws.Range("D19").Value = "11.47"
$ws.Range("E19").Value = "  +3.97%  "

# Row 20
$ws.Range("D20").Value = "367.25"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("D21").Value = "7.43"
$ws.Range("E21").Value = "  +0.93%  "

# Row 22
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("D23").Value = "4.86"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +4.74%  "

# Row 25
$ws.Range("D25").Value = "73.64"
$ws.Range("E25").Value = "  +0.68%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +1.41%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0000105"
$ws.Range("E28").Value = "  +6.39%  "

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.770.21"
$ws.Range("E29").Value = "  +1.27%  "

# Row 30
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31
$ws.Range("D31").Value = "575.04"
$ws.Range("E31").Value = "  -0.36%  "

# Row 32
$ws.Range("E32").Value = "  +5.18%  "

# Row 33
$ws.Range("D33").Value = "8.02"
$ws.Range("E33").Value = "  +4.75%  "

# Row 34
$ws.Range("E34").Value = "  +2.78%  "

# Row 35
$ws.Range("E35").Value = "  +3.51%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  +3.75%  "

# Row 38
$ws.Range("D38").Value = "160.40"
$ws.Range("E38").Value = "  +1.22%  "

# Row 39
$ws.Range("D39").Value = "19.32"
$ws.Range("E39").Value = "  +1.50%  "

# Row 40
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +3.69%  "

# Row 41
$ws.Range("E41").Value = "  +1.03%  "

# Row 42
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +3.52%  "

# Row 43
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  +3.85%  "

# Row 44
$ws.Range("D44").Value = "17.74"
$ws.Range("E44").Value = "  +3.62%  "

# Row 45
$ws.Range("E45").Value = "  +14.00%  "

# Row 46
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("D47").Value = "40.55"
$ws.Range("E47").Value = "  -0.28%  "

# Row 48
$ws.Range("D48").Value = "158.33"
$ws.Range("E48").Value = "  +3.40%  "

# Row 49
$ws.Range("E49").Value = "  +1.83%  "

# Row 50
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  +2.46%  "

# Row 51
$ws.Range("D51").Value = "21.99"
$ws.Range("E51").Value = "  +3.18%  "
